$wb = $excel.ActiveWorkbook

# 1 & 2: Fix reference names - these reference the same shared strings on both
#        "adipose" sheets (Vessel size + Vessel density), so update on both:
#    "Hul et al., 2012" -> "Van Hul et al., 2012"
#    "Lijnen, 2003" -> "Lijnen et al., 2003"
$wsVesselSizeAdipose = $wb.Worksheets.Item("Vessel size (adipose)")
$wsVesselSizeAdipose.Range("A7").Value = "Van Hul et al., 2012"
$wsVesselSizeAdipose.Range("A5").Value = "Lijnen et al., 2003"

$wsVesselDensityAdipose = $wb.Worksheets.Item("Vessel density (adipose)")
$wsVesselDensityAdipose.Range("A7").Value = "Van Hul et al., 2012"
$wsVesselDensityAdipose.Range("A5").Value = "Lijnen et al., 2003"

# 3: Delete the "Koyama, 2017" row (too large) from the "Vessel size (tumor)" table
$wsVesselSizeTumor = $wb.Worksheets.Item("Vessel size (tumor)")
$wsVesselSizeTumor.Rows.Item(9).Delete()

# Restore/update cursor position & active sheet to match the edited workbook state
$wsVesselSizeAdipose.Activate() | Out-Null
$wsVesselSizeAdipose.Range("A8").Select() | Out-Null

$wsVesselDensityAdipose.Activate() | Out-Null
$wsVesselDensityAdipose.Range("A5").Select() | Out-Null

$wsCbmRetina = $wb.Worksheets.Item("CBM (retina)")
$wsCbmRetina.Activate() | Out-Null
$wsCbmRetina.Range("B4").Select() | Out-Null

$wsCbmMuscle = $wb.Worksheets.Item("CBM (muscle)")
$wsCbmMuscle.Activate() | Out-Null
$wsCbmMuscle.Range("A16").Select() | Out-Null

$wsVesselSizeTumor.Activate() | Out-Null
$wsVesselSizeTumor.Range("B16").Select() | Out-Null
